$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Structural cleanup: drop the footnote hyperlink on sheet1!A8 (GLS row) ---
$ws1.Range("A8").Hyperlinks.Delete()

# --- Remove the "APV (LPV or LNAV/VNAV or RNP AR APCH)" row (row 9) from sheet1 ---
$ws1.Range("A9:E9").EntireRow.Delete()

# --- Drop the now-unused apt_nb / apt_share columns (D:E) from both sheets ---
$ws1.Range("D1:E1").EntireColumn.Delete()
$ws2.Range("D1:E1").EntireColumn.Delete()

# --- Update sheet1 (pbn_deployment) runway_nb / runway_share values ---
$ws1.Range("B2").Value = 1154
$ws1.Range("C2").Value = 0.71
$ws1.Range("B3").Value = 833
$ws1.Range("C3").Value = 0.51
$ws1.Range("B4").Value = 754
$ws1.Range("C4").Value = 0.46
$ws1.Range("B5").Value = 1173
$ws1.Range("C5").Value = 0.72
$ws1.Range("B6").Value = 36
$ws1.Range("C6").Value = 0.02
$ws1.Range("B7").Value = 775
$ws1.Range("C7").Value = 0.47
$ws1.Range("B8").Value = 40
$ws1.Range("C8").Value = 0.02
$ws1.Range("B9").Value = 948
$ws1.Range("C9").Value = 0.58

# --- Update sheet2 (ils_deployment) runway_nb / runway_share values ---
$ws2.Range("B2").Value = 512
$ws2.Range("C2").Value = 0.313
$ws2.Range("B3").Value = 263
$ws2.Range("C3").Value = 0.16

# --- Defined name _ftn1 now points one row higher since the APV row was removed ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "pbn_deployment!_ftn1") {
        $n.RefersTo = "=pbn_deployment!`$A`$12"
    }
}

# --- Restore view state: selections on each sheet, sheet2 as the active tab ---
$ws1.Activate()
$ws1.Range("G12").Select()
$ws2.Activate()
$ws2.Range("I13").Select()

Write-Output "edit complete"
